# Update countries & provincias Spain
# Refreshes the "Pais" COVID-19 snapshot table: a handful of countries swapped
# rank (so their row shows a different country name) and several rows' stats
# were refreshed, plus the "last updated" timestamp advanced from 01:31 to 02:48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 02:48"

# --- Country name changes (rows whose rank now belongs to a different
#     country because of the data refresh) --------------------------------
$ws.Range("A114").Value = "Nicaragua"
$ws.Range("A115").Value = "Suazilandia"

$ws.Range("A117").Value = "Surinam"
$ws.Range("A118").Value = "Cuba"
$ws.Range("A119").Value = "Mozambique"

$ws.Range("A142").Value = "Aruba"
$ws.Range("A143").Value = "Jordania"
$ws.Range("A144").Value = "Sierra Leona"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ------------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6257171
$ws.Range("C4").Value = 41579
$ws.Range("D4").Value = 3485908
$ws.Range("E4").Value = 2582389
$ws.Range("G4").Value = 1138
$ws.Range("H4").Value = 188874

# Indonesia (row 36)
$ws.Range("B36").Value = 93552
$ws.Range("C36").Value = 570
$ws.Range("D36").Value = 66974
$ws.Range("E36").Value = 24560
$ws.Range("G36").Value = 16
$ws.Range("H36").Value = 2018

# Nicaragua (row 114)
$ws.Range("B114").Value = 4668
$ws.Range("C114").Value = 174
$ws.Range("D114").Value = 2913
$ws.Range("E114").Value = 1614
$ws.Range("G114").Value = 4
$ws.Range("H114").Value = 141

# Suazilandia (row 115)
$ws.Range("B115").Value = 4618
$ws.Range("C115").Value = 41
$ws.Range("D115").Value = 3562
$ws.Range("E115").Value = 962
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 94

# Surinam (row 117)
$ws.Range("B117").Value = 4089
$ws.Range("C117").Value = 55
$ws.Range("D117").Value = 3171
$ws.Range("E117").Value = 846
$ws.Range("H117").Value = 72

# Cuba (row 118)
$ws.Range("B118").Value = 4065
$ws.Range("C118").Value = 33
$ws.Range("D118").Value = 3395
$ws.Range("E118").Value = 575
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 95

# Mozambique (row 119)
$ws.Range("B119").Value = 4039
$ws.Range("C119").Value = 123
$ws.Range("D119").Value = 2170
$ws.Range("E119").Value = 1846
$ws.Range("H119").Value = 23

# Tunez (row 123, name unchanged - stats only)
$ws.Range("B123").Value = 3963
$ws.Range("C123").Value = 160
$ws.Range("D123").Value = 1624
$ws.Range("E123").Value = 2259
$ws.Range("G123").Value = 3
$ws.Range("H123").Value = 80

# row 138
$ws.Range("B138").Value = 2276
$ws.Range("C138").Value = 59
$ws.Range("E138").Value = 1425

# Aruba (row 142)
$ws.Range("B142").Value = 2104
$ws.Range("C142").Value = 98
$ws.Range("D142").Value = 857
$ws.Range("E142").Value = 1237
$ws.Range("H142").Value = 10

# Jordania (row 143)
$ws.Range("B143").Value = 2097
$ws.Range("C143").Value = 63
$ws.Range("D143").Value = 1564
$ws.Range("E143").Value = 518
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 15

# Sierra Leona (row 144)
$ws.Range("B144").Value = 2028
$ws.Range("C144").Value = 6
$ws.Range("D144").Value = 1594
$ws.Range("E144").Value = 363
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 71

# row 167
$ws.Range("D167").Value = 853
$ws.Range("E167").Value = 28

# row 188
$ws.Range("D188").Value = 204
$ws.Range("E188").Value = 0

# row 189
$ws.Range("B189").Value = 176
$ws.Range("C189").Value = 2
$ws.Range("E189").Value = 23

# row 197
$ws.Range("B197").Value = 61
$ws.Range("C197").Value = 1
$ws.Range("E197").Value = 4

# row 212
$ws.Range("B212").Value = 16
$ws.Range("C212").Value = 1
$ws.Range("E212").Value = 9
